$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-09-07 12:39:38"

for ($r = 2; $r -le 7; $r++) {
    $ws.Cells.Item($r, 1).Value = $newTimestamp
}
